# Refresh the 'cryptos' price/volume snapshot (GitHub Actions scheduled update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 35 & 36 swap ranking order: Kaspa moves up to rank 34 (row 35),
#     Monero drops to rank 35 (row 36). Update Coin / Link / Price / Volume together.
$ws.Range("B35").Value = 'Kaspa'
$ws.Range("C35").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.121'
$ws.Range("E35").Value = '  +2.76%  '

$ws.Range("B36").Value = 'Monero'
$ws.Range("C36").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '162.94'
$ws.Range("E36").Value = '  +1.62%  '

# --- Remaining Price / Volume(1h) corrections ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '68.026.82'
$ws.Range("E2").Value = '  +1.39%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.508.52'
$ws.Range("E3").Value = '  +1.09%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '589.14'
$ws.Range("E5").Value = '  +0.96%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '176.84'
$ws.Range("E6").Value = '  +3.54%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("E8").Value = '  +0.71%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.144'
$ws.Range("E9").Value = '  +5.01%  '
$ws.Range("E10").Value = '  -0.63%  '
$ws.Range("E11").Value = '  +1.95%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.97'
$ws.Range("E12").Value = '  +0.87%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.979.02'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.75'
$ws.Range("E14").Value = '  +1.35%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '67.896.54'
$ws.Range("E15").Value = '  +1.39%  '
$ws.Range("E16").Value = '  +1.18%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.568.85'
$ws.Range("E17").Value = '  -1.08%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '10.99'
$ws.Range("E18").Value = '  -0.17%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.54'
$ws.Range("E19").Value = '  +1.37%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '352.81'
$ws.Range("E20").Value = '  +1.25%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.13'
$ws.Range("E21").Value = '  +2.36%  '
$ws.Range("E22").Value = '  +0.09%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '71.03'
$ws.Range("E23").Value = '  +3.76%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.28'
$ws.Range("E24").Value = '  +0.99%  '
$ws.Range("E25").Value = '  -2.60%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.638.61'
$ws.Range("E27").Value = '  +0.94%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.997'
$ws.Range("E28").Value = '  -0.13%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0₃0917'
$ws.Range("E29").Value = '  +1.18%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '509.79'
$ws.Range("E30").Value = '  -0.26%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.80'
$ws.Range("E31").Value = '  +0.93%  '
$ws.Range("E32").Value = '  +2.85%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.77'
$ws.Range("E33").Value = '  +0.66%  '
$ws.Range("E34").Value = '  +0.01%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '18.42'
$ws.Range("E37").Value = '  +1.02%  '
$ws.Range("E38").Value = '  -0.29%  '
$ws.Range("E39").Value = '  +0.08%  '
$ws.Range("E40").Value = '  +0.02%  '
$ws.Range("E41").Value = '  +2.97%  '
$ws.Range("E42").Value = '  +0.50%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.87'
$ws.Range("E43").Value = '  +1.70%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.49'
$ws.Range("E44").Value = '  +5.62%  '
$ws.Range("E45").Value = '  +3.50%  '
$ws.Range("E46").Value = '  +2.80%  '
$ws.Range("E47").Value = '  +4.41%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.520'
$ws.Range("E48").Value = '  +1.14%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0743'
$ws.Range("E49").Value = '  +2.05%  '
$ws.Range("E50").Value = '  +1.80%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.586'
$ws.Range("E51").Value = '  +0.52%  '
